# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.827.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.735.76"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5164"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.39"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06110"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.736.77"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07036"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.10"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6329"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.490"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.72"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.830.91"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006631"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.954.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.125"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.708"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.119"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.508"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.768"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "101.65"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08264"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.669"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.468"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04478"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9722"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6117"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.656"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01573"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.925"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.18"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3801"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.987"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7188"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05370"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1117"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.181"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.85"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.560"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.86%  "
